$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.642.08"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "1.844.91"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'259.87"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D7").Value = "'0.5280"
$ws.Range("E7").Value = "  +1.76%  "
$ws.Range("D8").Value = "'0.3157"
$ws.Range("E8").Value = "  -3.26%  "
$ws.Range("D9").Value = "'0.06802"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "'18.97"
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("D11").Value = "'0.7858"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").Value = "'0.07784"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").Value = "1.857.65"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "'88.31"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "'5.017"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "'0.000007924"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "26.652.72"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "2.088.58"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'5.988"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").Value = "'9.351"
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("D25").Value = "'2.229"
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("D26").Value = "'143.00"
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").Value = "'17.04"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("D29").Value = "'111.05"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("D31").Value = "'0.08714"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").Value = "'0.04886"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").Value = "'0.7323"
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("D35").Value = "'1.142"
$ws.Range("E35").Value = "  +1.59%  "
$ws.Range("D38").Value = "'2.297"
$ws.Range("E38").Value = "  +3.76%  "
$ws.Range("D39").Value = "'0.01733"
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("D40").Value = "'0.4818"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").Value = "'0.9010"
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("D42").Value = "'109.84"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").Value = "'5.939"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "'7.713"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'0.4202"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("D47").Value = "'9.124"
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("D48").Value = "'0.1244"
$ws.Range("E48").Value = "  +1.66%  "
$ws.Range("D49").Value = "'0.05829"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").Value = "'0.8960"
$ws.Range("E51").Value = "  +1.25%  "
